$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 1).Value = 1013
$ws.Cells.Item(19, 2).Value = "ok test"
$ws.Cells.Item(19, 3).Value = 1522
$ws.Cells.Item(19, 4).Value = "PRJ-05"

# Row 20
$ws.Cells.Item(20, 1).Value = 1014
$ws.Cells.Item(20, 2).Value = "test PRJ"
$ws.Cells.Item(20, 3).Value = 15
$ws.Cells.Item(20, 4).Value = "PRJ-564"

$ws.Range("D20").Select()
